# Update "想去人数" (number of people interested) values that changed
# between two scrapes of the 合肥-漫展信息 data, as published to gh-pages.
#
# Sheet "展览"   (sheet1): F5 4954->4962, F9 751->752, F11 3->4
# Sheet "演出"   (sheet2): F2 28->29
# Sheet "本地生活" (sheet3): no changes
# Sheet "全部类型" (sheet4): F5 4954->4962, F9 751->752, F10 28->29, F12 3->4

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 4962
$wsExhibition.Range("F9").Value = 752
$wsExhibition.Range("F11").Value = 4

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 29

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 4962
$wsAll.Range("F9").Value = 752
$wsAll.Range("F10").Value = 29
$wsAll.Range("F12").Value = 4
